# Cambio de política de 32 a 16 estados en cada carril
#
# 1) Rename the existing sheet "Sheet1" -> "32"
# 2) Insert a new worksheet "16" right after it, and make it the active sheet
# 3) Fill it with the 16-state policy table (A2:D17), in the big 22pt Arial
#    "poster" style used for this kind of summary sheet
# 4) Match column widths / row heights / view (zoom, selection) / print
#    margins as closely as the object model allows

$wb = $excel.ActiveWorkbook

# --- 1) rename the original sheet ------------------------------------------------
$ws32 = $wb.Worksheets.Item(1)
$ws32.Name = "32"

# --- 2) add the new sheet right after "32" and rename it --------------------------
$ws16 = $wb.Worksheets.Add([System.Type]::Missing, $ws32)
$ws16.Name = "16"

# --- 3) column widths (set before writing values / fonts) -------------------------
$ws16.Columns.Item(1).ColumnWidth = 37.333333333333336   # -> stored width ~38.16
$ws16.Columns.Item(2).ColumnWidth = 29.333333333333332   # -> stored width ~30.13
$ws16.Columns.Item(3).ColumnWidth = 28.5                 # -> stored width ~29.35
$ws16.Columns.Item(4).ColumnWidth = 28.333333333333332   # -> stored width ~29.16

# --- 4) the 16-state policy table --------------------------------------------------
$rows = @(
  @("Policy(free_N(0)=0", " free_NW(0)=0", " free_SW(0)=0", " free_W(0)=0) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=0", " free_SW(0)=0", " free_W(0)=0) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=1", " free_SW(0)=0", " free_W(0)=0) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=1", " free_SW(0)=0", " free_W(0)=0) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=0", " free_SW(0)=1", " free_W(0)=0) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=0", " free_SW(0)=1", " free_W(0)=0) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=1", " free_SW(0)=1", " free_W(0)=0) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=1", " free_SW(0)=1", " free_W(0)=0) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=0", " free_SW(0)=0", " free_W(0)=1) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=0", " free_SW(0)=0", " free_W(0)=1) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=1", " free_SW(0)=0", " free_W(0)=1) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=1", " free_SW(0)=0", " free_W(0)=1) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=0", " free_SW(0)=1", " free_W(0)=1) = keep_distance"),
  @("Policy(free_N(0)=1", " free_NW(0)=0", " free_SW(0)=1", " free_W(0)=1) = cruise"),
  @("Policy(free_N(0)=0", " free_NW(0)=1", " free_SW(0)=1", " free_W(0)=1) = change_lane"),
  @("Policy(free_N(0)=1", " free_NW(0)=1", " free_SW(0)=1", " free_W(0)=1) = cruise")
)

$r = 2
foreach ($row in $rows) {
    $ws16.Cells.Item($r, 1).Value = $row[0]
    $ws16.Cells.Item($r, 2).Value = $row[1]
    $ws16.Cells.Item($r, 3).Value = $row[2]
    $ws16.Cells.Item($r, 4).Value = $row[3]
    $ws16.Rows.Item($r).RowHeight = 26.8
    $r++
}

# --- 5) big poster font across the whole table -------------------------------------
$tbl = $ws16.Range("A2:D17")
$tbl.Font.Name = "Arial"
$tbl.Font.Size = 22

# --- 6) view: zoom 71%, gridlines/headers on, selection at G8 ----------------------
$ws16.Activate()
$ws16.Range("G8").Select()
$excel.ActiveWindow.Zoom = 71
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# --- 7) print setup: margins / orientation / headers & footers ---------------------
$ws16.PageSetup.LeftMargin = 0.7875 * 72
$ws16.PageSetup.RightMargin = 0.7875 * 72
$ws16.PageSetup.TopMargin = 1.05277777777778 * 72
$ws16.PageSetup.BottomMargin = 1.05277777777778 * 72
$ws16.PageSetup.HeaderMargin = 0.7875 * 72
$ws16.PageSetup.FooterMargin = 0.7875 * 72
$ws16.PageSetup.Orientation = 1
$ws16.PageSetup.PaperSize = 1
$ws16.PageSetup.PrintHeadings = $false
$ws16.PageSetup.PrintGridlines = $false
$ws16.PageSetup.CenterHorizontally = $false
$ws16.PageSetup.CenterVertically = $false
$ws16.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws16.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

Write-Host "done"
